$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (D) and 1h volume-change (E) updates.
# Price cells that look like plain decimals must be forced to text
# (NumberFormat "@") before assignment so Excel does not coerce them
# into floating-point numbers; the style is reset to Normal afterwards
# so no stray number-format style is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.677.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.368.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.76%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.368.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.14%  "
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.799.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.602.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.367.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.22%  "
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "570.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.487.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("E30").Value = "  -4.05%  "
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("E34").Value = "  -5.49%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  -5.52%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "146.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("E45").Value = "  -5.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0281"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +19.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.584"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "
